$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Wednesday hours for the week commencing 43178 (row 10) from 0 to 1.5
$ws.Range("D10").Value = 1.5

# Update the active selection on the sheet to match the saved view state
$ws.Range("H28").Select()

$wb.Save()
